$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (B2, C2) ---
$ws.Range("B2").Value = "sai rồi nè 3"
$ws.Range("C2").Value = "Sai rồi"

# --- Add new row 3 ---
$ws.Range("A3").Value = "Phân bón fail fail 2 "
$ws.Range("B3").Value = "Công ty Adama"
$ws.Range("C3").Value = "sai nữa nè 3"
$ws.Range("D3").Value = "2L"
$ws.Range("E3").Value = 50000
$ws.Range("F3").Value = 18000
$ws.Range("G3").Value = 5

# --- Add new row 4 ---
$ws.Range("A4").Value = "Phân bón fail fail"
$ws.Range("B4").Value = "sai rồi nè 5"
$ws.Range("C4").Value = "Sai rồi"
$ws.Range("D4").Value = "2L"
$ws.Range("E4").Value = 50000
$ws.Range("F4").Value = 18000
$ws.Range("G4").Value = 5

# --- Add new row 5 ---
$ws.Range("A5").Value = "Phân bón fail fail 2 "
$ws.Range("B5").Value = "Công ty Adama"
$ws.Range("C5").Value = "sai nữa nè 5"
$ws.Range("D5").Value = "2L"
$ws.Range("E5").Value = 50000
$ws.Range("F5").Value = 18000
$ws.Range("G5").Value = 5

# --- Copy formatting from row 2 template cells so new rows match style pattern ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null

$ws.Range("E2:G2").Copy() | Out-Null
$ws.Range("E4:G4").PasteSpecial(-4122) | Out-Null

$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A3:G3").PasteSpecial(-4122) | Out-Null
$ws.Range("A5:G5").PasteSpecial(-4122) | Out-Null

# fix B/C styles on rows 3 & 5 to use the "alt" look (style s=8 / s=7)
$ws.Range("B3").Style = $ws.Range("C2").Style
$ws.Range("C3").Style = $ws.Range("B2").Style
$ws.Range("B5").Style = $ws.Range("C2").Style
$ws.Range("C5").Style = $ws.Range("B2").Style

$excel.CutCopyMode = 0

$ws.Range("C5").Select() | Out-Null
